# Updates crypto price-tracker sheet (cryptos.xlsx) for the latest GitHub
# Actions data-refresh run: Price (D), Volume(1h) (E) and Hora (G) columns
# for rows 2-51. Values are written as TEXT (not numbers/percentages) to
# match the source data's inline-string cell type, so each cell's
# NumberFormat is forced to "@" (Text) immediately before the assignment -
# otherwise Excel's normal type-inference would turn a numeric-looking
# string like "303.19" into a number, or "-0.63%" into a percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 2; D = "303.19"; E = "-0.63%"; G = "13" }
    @{ Row = 3; D = "37.35"; E = "6.57%"; G = "13" }
    @{ Row = 4; D = "5.006"; E = "-3.27%"; G = "13" }
    @{ Row = 5; D = "0.07873"; E = "0.54%"; G = "13" }
    @{ Row = 6; D = "2.211"; E = "-4.12%"; G = "13" }
    @{ Row = 7; D = $null; E = "-0.59%"; G = "13" }
    @{ Row = 8; D = "4.023"; E = "0.93%"; G = "13" }
    @{ Row = 9; D = "0.9194"; E = "-0.70%"; G = "13" }
    @{ Row = 10; D = "0.09549"; E = "-4.41%"; G = "13" }
    @{ Row = 11; D = "0.1881"; E = "2.65%"; G = "13" }
    @{ Row = 12; D = "0.08540"; E = "-0.17%"; G = "13" }
    @{ Row = 13; D = "0.03603"; E = "7.19%"; G = "13" }
    @{ Row = 14; D = "0.09968"; E = "0.59%"; G = "13" }
    @{ Row = 15; D = "0.001478"; E = "0.02%"; G = "13" }
    @{ Row = 16; D = "0.005680"; E = "-1.40%"; G = "13" }
    @{ Row = 17; D = "3.461"; E = "-0.70%"; G = "13" }
    @{ Row = 18; D = "2.249"; E = "5.75%"; G = "13" }
    @{ Row = 19; D = "0.3419"; E = "-0.68%"; G = "13" }
    @{ Row = 20; D = $null; E = "-0.57%"; G = "13" }
    @{ Row = 21; D = "4.746"; E = "4.21%"; G = "13" }
    @{ Row = 22; D = "0.2203"; E = "-7.63%"; G = "13" }
    @{ Row = 23; D = "0.04582"; E = "-1.37%"; G = "13" }
    @{ Row = 24; D = "0.001229"; E = "0.55%"; G = "13" }
    @{ Row = 25; D = "0.004786"; E = "7.60%"; G = "13" }
    @{ Row = 26; D = "0.0001398"; E = "7.74%"; G = "13" }
    @{ Row = 27; D = $null; E = "39.96%"; G = "13" }
    @{ Row = 28; D = $null; E = $null; G = "13" }
    @{ Row = 29; D = $null; E = $null; G = "13" }
    @{ Row = 30; D = $null; E = $null; G = "13" }
    @{ Row = 31; D = $null; E = $null; G = "13" }
    @{ Row = 32; D = $null; E = $null; G = "13" }
    @{ Row = 33; D = $null; E = $null; G = "13" }
    @{ Row = 34; D = $null; E = $null; G = "13" }
    @{ Row = 35; D = $null; E = $null; G = "13" }
    @{ Row = 36; D = $null; E = $null; G = "13" }
    @{ Row = 37; D = $null; E = $null; G = "13" }
    @{ Row = 38; D = $null; E = $null; G = "13" }
    @{ Row = 39; D = "0.01806"; E = "2.71%"; G = "13" }
    @{ Row = 40; D = $null; E = "-0.46%"; G = "13" }
    @{ Row = 41; D = "0.008155"; E = "6.03%"; G = "13" }
    @{ Row = 42; D = "0.1395"; E = "-1.16%"; G = "13" }
    @{ Row = 43; D = "0.007557"; E = "7.04%"; G = "13" }
    @{ Row = 44; D = "0.002207"; E = "0.06%"; G = "13" }
    @{ Row = 45; D = "0.01042"; E = "3.92%"; G = "13" }
    @{ Row = 46; D = "0.00006153"; E = "2.74%"; G = "13" }
    @{ Row = 47; D = "0.00000000749"; E = "0.05%"; G = "13" }
    @{ Row = 48; D = "0.0005800"; E = "-0.01%"; G = "13" }
    @{ Row = 49; D = "7.114"; E = "22.72%"; G = "13" }
    @{ Row = 50; D = $null; E = "0.18%"; G = "13" }
    @{ Row = 51; D = "0.00002098"; E = "0.05%"; G = "13" }
)

foreach ($u in $rowUpdates) {
    foreach ($col in @("D", "E", "G")) {
        $newVal = $u[$col]
        if ($null -eq $newVal) { continue }
        $cell = $ws.Range("$col$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
    }
}
